$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (id, speaker_variant) pairs for rows 2-21.
# The "id" (column B) is recomputed as a slug of the speaker_variant (column C):
# lower-cased, spaces replaced with hyphens.
$rows = @(
    @("#belje", "Belje"),
    @("#schout", "Schout"),
    @("#kniertie", "Kniertie"),
    @("#kniertje", "Kniertje"),
    @("#vólkert", "Vólkert"),
    @("#volkert", "Volkert"),
    @("#kristoffel", "Kristoffel"),
    @("#filibeet", "Filibeet"),
    @("#geertruy", "Geertruy"),
    @("#filibert", "Filibert"),
    @("#dienaar", "Dienaar"),
    @("#geertrui", "Geertrui"),
    @("#ioris", "Ioris"),
    @("#1-dienaar", "1 Dienaar"),
    @("#belt", "Belt"),
    @("#joris", "Joris"),
    @("#agata", "Agata"),
    @("#filibert-het-stukken-scheurende", "FILIBERT het stukken scheurende"),
    @("#dirk", "Dirk"),
    @("#bely", "Bely")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $pair = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $pair[0]
    $ws.Cells.Item($r, 3).Value = $pair[1]
    # is_prefered (column D) is cleared for every data row in the new export.
    $ws.Cells.Item($r, 4).Value = ""
}
